$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the header date row: " Wednesday"/01/04/2020 -> " Tuesday"/07/04/2020
$ws.Range("A1").Value = "  Tuesday"

# B1 holds the date as literal text (not a real date serial). Assigning the
# string straight to .Value would have Excel auto-convert it to a date
# serial, so instead enter it as a quoted text formula, then flatten that
# formula down to its literal value/shared-string in place. This keeps B1 on
# its original (General) number format/style instead of picking up a new one.
$c = $ws.Range("B1")
$c.Formula = "=""07/04/2020"""
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues

# Row 3: A3 placeholder name becomes a separator line
$ws.Range("A3").Value = "___________"

# Row 4: A4 placeholder name is removed (cell left blank, like the rows below it)
$ws.Range("A4").Value = ""

# Roster shrinks by two (Shakti & Dganit leave) - the remaining names shuffle up
$ws.Range("B8").Value = "Supriti"
$ws.Range("B9").Value = "Anna"
